# Commit: "Update database to include EV loads"
#
# Adds a new "Ev_kW" column (EV charging loads) to the INTERNAL_LOADS
# sheet of the use-types properties workbook, with a default value of 0
# for every occupancy-type row, and leaves the INTERNAL_LOADS sheet/cell
# N1 as the active selection (mirroring the author re-saving the file
# right after typing the new header).

$wb  = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item("INTERNAL_LOADS")

# Existing data occupies columns A:M, rows 1 (header) through 20 (data).
# The new column goes right after M, i.e. column N (14).
$lastRow = 20
$srcCol  = "M"
$dstCol  = "N"

# --- Header cell (N1): label + copy the header's formatting from M1 ---
$ws.Range($dstCol + "1").Value = "Ev_kW"
$ws.Range($srcCol + "1").Copy() | Out-Null
$ws.Range($dstCol + "1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- Data cells (N2:N20): default to 0, formatted like the M column ---
$ws.Range($srcCol + "2:" + $srcCol + $lastRow).Copy() | Out-Null
$ws.Range($dstCol + "2:" + $dstCol + $lastRow).PasteSpecial(-4122) | Out-Null

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 14).Value = 0
}

$excel.CutCopyMode = 0

# --- Restore/ set the focused view: INTERNAL_LOADS active, N1 selected ---
$ws.Activate()
$ws.Range("N1").Select() | Out-Null

Write-Host "Added Ev_kW column to INTERNAL_LOADS"
